$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Remove the old "_GoBack" bookmark that sits right after
#    "SE TENDRÁ POR NO PUESTA. " (it gets relocated further down).
# ------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# ------------------------------------------------------------------
# 2) Rework the pagaré amount clause:
#    "(${#contrato}{capitalSolicitado}{/contrato}) ___________ PESOS
#    CERO CENTAVOS MONEDA NACIONAL, VALOR RECIBIDO..."
#    becomes
#    "(${#contrato}{capitalSolicitado}{/contrato}), {letra} VALOR
#    RECIBIDO..." with a fresh "_GoBack" bookmark right before "VALOR".
# ------------------------------------------------------------------

# 2a) Drop the spelled-out placeholder amount and keep just "),".
$rngAmount = $d.Content.Duplicate
$rngAmount.Find.Execute(") ___________ PESOS CERO CENTAVOS MONEDA NACIONAL,", `
    $false, $false, $false, $false, $false, $true, 1, $false, "),", 2) | Out-Null

# 2b) Insert " {letra}" (bold) right after the "),".
$rngLetra = $d.Content.Duplicate
$rngLetra.Find.Execute("{/contrato}),") | Out-Null
$rngLetra.Collapse(0)
$rngLetra.InsertAfter(" {letra}")
$rngLetra.Bold = 1

# 2c) Re-insert the "_GoBack" bookmark immediately before "VALOR RECIBIDO".
$rngValor = $d.Content.Duplicate
$rngValor.Find.Execute("VALOR RECIBIDO A MI ENTERA SATISFACCIÓN") | Out-Null
$rngMark = $d.Range($rngValor.Start, $rngValor.Start)
$rngMark.InsertBefore("Z")
$d.Bookmarks.Add("_GoBack", $rngMark)
$rngMark.Text = ""
